$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 932.6667
$ws.Range("I4").Value = 932.6667
$ws.Range("K4").Value = 932.6667
$ws.Range("M4").Value = -818.6667
$ws.Range("H8").Value = 1015.4
$ws.Range("I8").Value = 1015.4
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 3046.2
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -2907.2
$ws.Range("N8").ClearContents()
$ws.Range("H9").Value = 125.21429
$ws.Range("I9").Value = 107
$ws.Range("J9").Value = 170.75
$ws.Range("K9").Value = 107
$ws.Range("L9").Value = 170.75
$ws.Range("M9").Value = 62
$ws.Range("N9").Value = -508.75
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H17").Value = 319.21622
$ws.Range("J17").Value = 321.97223
$ws.Range("L17").Value = 965.91669
$ws.Range("N17").Value = -1301.91669
$ws.Range("H18").Value = 142857630
$ws.Range("I18").Value = 272.25
$ws.Range("J18").Value = 333334100
$ws.Range("K18").Value = 272.25
$ws.Range("L18").Value = 333334100
$ws.Range("M18").Value = 11.75
$ws.Range("N18").Value = -333334668
$ws.Range("H19").Value = 639.625
$ws.Range("I19").Value = 614
$ws.Range("J19").Value = 659.55554
$ws.Range("K19").Value = 614
$ws.Range("L19").Value = 659.55554
$ws.Range("M19").Value = -439
$ws.Range("N19").Value = -1009.55554
$ws.Range("H38").Value = 2575.6667
$ws.Range("I38").Value = 2128
$ws.Range("J38").Value = 7500
$ws.Range("K38").Value = 6384
$ws.Range("L38").Value = 22500
$ws.Range("M38").Value = -6012
$ws.Range("N38").Value = -23244
$ws.Range("H40").Value = 4519
$ws.Range("I40").Value = 3255.2856
$ws.Range("J40").Value = 5624.75
$ws.Range("K40").Value = 3255.2856
$ws.Range("L40").Value = 5624.75
$ws.Range("M40").Value = -3080.2856
$ws.Range("N40").Value = -5974.75
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("H106").Value = 1938.8334
$ws.Range("I106").Value = 1592.6154
$ws.Range("J106").Value = 2839
$ws.Range("K106").Value = 1592.6154
$ws.Range("L106").Value = 2839
$ws.Range("M106").Value = -961.6153999999999
$ws.Range("N106").Value = -4101
$ws.Range("H107").Value = 920
$ws.Range("I107").Value = 920
$ws.Range("K107").Value = 920
$ws.Range("M107").Value = 1000
$ws.Range("H116").Value = 6027
$ws.Range("I116").Value = 10994
$ws.Range("J116").Value = 5317.4287
$ws.Range("K116").Value = 10994
$ws.Range("L116").Value = 5317.4287
$ws.Range("M116").Value = -7552
$ws.Range("N116").Value = -12201.4287
$ws.Range("H132").Value = 7420.8335
$ws.Range("I132").Value = 6732.0454
$ws.Range("K132").Value = 20196.1362
$ws.Range("M132").Value = -17666.1362
$ws.Range("H137").Value = 4122.75
$ws.Range("I137").Value = 1632.3334
$ws.Range("J137").Value = 6875.316
$ws.Range("K137").Value = 4897.0002
$ws.Range("L137").Value = 20625.948
$ws.Range("M137").Value = -2347.0002
$ws.Range("N137").Value = -25725.948
$ws.Range("H138").Value = 4231.5
$ws.Range("J138").Value = 4637.807
$ws.Range("L138").Value = 13913.421
$ws.Range("N138").Value = -24193.421
$ws.Range("H141").Value = 3648.25
$ws.Range("I141").Value = 3648.25
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 10944.75
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -5764.75
$ws.Range("N141").ClearContents()

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 72.375
$ws.Range("I5").Value = 71.28570999999999
$ws.Range("K5").Value = 71.28570999999999
$ws.Range("M5").Value = 40.71429000000001
$ws.Range("H32").Value = 2018074
$ws.Range("I32").Value = 1003627.4
$ws.Range("K32").Value = 1003627.4
$ws.Range("M32").Value = -1003340.4
$ws.Range("H61").Value = 6383.25
$ws.Range("I61").Value = 2198.125
$ws.Range("J61").Value = 14753.5
$ws.Range("K61").Value = 2198.125
$ws.Range("L61").Value = 14753.5
$ws.Range("M61").Value = -1986.125
$ws.Range("N61").Value = -15177.5
$ws.Range("H74").Value = 23260952
$ws.Range("I74").Value = 4314.5713
$ws.Range("J74").Value = 34488296
$ws.Range("K74").Value = 4314.5713
$ws.Range("L74").Value = 34488296
$ws.Range("M74").Value = -3440.5713
$ws.Range("N74").Value = -34490044
$ws.Range("H77").Value = 23260952
$ws.Range("I77").Value = 4314.5713
$ws.Range("J77").Value = 34488296
$ws.Range("K77").Value = 21572.8565
$ws.Range("L77").Value = 172441480
$ws.Range("M77").Value = -17204.8565
$ws.Range("N77").Value = -172450216
$ws.Range("H122").Value = 3962.45
$ws.Range("I122").Value = 3504.4546
$ws.Range("J122").Value = 4522.222
$ws.Range("K122").Value = 10513.3638
$ws.Range("L122").Value = 13566.666
$ws.Range("M122").Value = -8063.363799999999
$ws.Range("N122").Value = -18466.666
$ws.Range("H132").Value = 4531163.5
$ws.Range("I132").Value = 6997798.5
$ws.Range("K132").Value = 20993395.5
$ws.Range("M132").Value = -20990865.5
$ws.Range("H136").Value = 6383.25
$ws.Range("I136").Value = 2198.125
$ws.Range("J136").Value = 14753.5
$ws.Range("K136").Value = 6594.375
$ws.Range("L136").Value = 44260.5
$ws.Range("M136").Value = -4044.375
$ws.Range("N136").Value = -49360.5

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 72.375
$ws.Range("I4").Value = 71.28570999999999
$ws.Range("K4").Value = 71.28570999999999
$ws.Range("M4").Value = 43.71429000000001
$ws.Range("H20").Value = 29766264
$ws.Range("I20").Value = 34726996
$ws.Range("J20").Value = 1865.25
$ws.Range("K20").Value = 34726996
$ws.Range("L20").Value = 1865.25
$ws.Range("M20").Value = -34726749
$ws.Range("N20").Value = -2359.25
$ws.Range("H35").Value = 54037
$ws.Range("J35").Value = 54037
$ws.Range("L35").Value = 54037
$ws.Range("N35").Value = -54657
$ws.Range("H80").Value = 76923400
$ws.Range("J80").Value = 407.125
$ws.Range("L80").Value = 407.125
$ws.Range("N80").Value = -2403.125
$ws.Range("H83").Value = 76923400
$ws.Range("J83").Value = 407.125
$ws.Range("L83").Value = 2035.625
$ws.Range("N83").Value = -12019.625
$ws.Range("H99").Value = 3394.6667
$ws.Range("I99").Value = 3020
$ws.Range("K99").Value = 3020
$ws.Range("M99").Value = -1522
$ws.Range("H105").Value = 11819629
$ws.Range("I105").Value = 556546.4
$ws.Range("K105").Value = 556546.4
$ws.Range("M105").Value = -554799.4
$ws.Range("H134").Value = 2139.25
$ws.Range("I134").Value = 1668.5834
$ws.Range("J134").Value = 4963.25
$ws.Range("K134").Value = 5005.7502
$ws.Range("L134").Value = 14889.75
$ws.Range("M134").Value = -2470.7502
$ws.Range("N134").Value = -19959.75
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 328.66666
$ws.Range("I7").Value = 342.7143
$ws.Range("J7").Value = 279.5
$ws.Range("K7").Value = 342.7143
$ws.Range("L7").Value = 279.5
$ws.Range("M7").Value = -229.7143
$ws.Range("N7").Value = -505.5
$ws.Range("H16").Value = 1500
$ws.Range("I16").Value = 1500
$ws.Range("K16").Value = 1500
$ws.Range("M16").Value = -1213
$ws.Range("H31").Value = 2606725.5
$ws.Range("I31").Value = 2028.75
$ws.Range("J31").Value = 3474957.8
$ws.Range("K31").Value = 2028.75
$ws.Range("L31").Value = 3474957.8
$ws.Range("M31").Value = -1733.75
$ws.Range("N31").Value = -3475547.8
$ws.Range("H34").Value = 2606725.5
$ws.Range("I34").Value = 2028.75
$ws.Range("J34").Value = 3474957.8
$ws.Range("K34").Value = 2028.75
$ws.Range("L34").Value = 3474957.8
$ws.Range("M34").Value = -1826.75
$ws.Range("N34").Value = -3475361.8
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H58").Value = 3728.3635
$ws.Range("I58").Value = 2666.6667
$ws.Range("K58").Value = 2666.6667
$ws.Range("M58").Value = -2463.6667
$ws.Range("H113").Value = 1500
$ws.Range("I113").Value = 1500
$ws.Range("K113").Value = 1500
$ws.Range("M113").Value = 670
$ws.Range("H122").Value = 1484.5714
$ws.Range("I122").Value = 1713.6666
$ws.Range("J122").Value = 1179.1111
$ws.Range("K122").Value = 5140.9998
$ws.Range("L122").Value = 3537.3333
$ws.Range("M122").Value = -2690.9998
$ws.Range("N122").Value = -8437.3333
$ws.Range("H132").Value = 3428.72
$ws.Range("I132").Value = 3377.2632
$ws.Range("J132").Value = 3591.6667
$ws.Range("K132").Value = 10131.7896
$ws.Range("L132").Value = 10775.0001
$ws.Range("M132").Value = -7601.7896
$ws.Range("N132").Value = -15835.0001
$ws.Range("H136").Value = 3728.3635
$ws.Range("I136").Value = 2666.6667
$ws.Range("K136").Value = 8000.000100000001
$ws.Range("M136").Value = -5450.000100000001

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 1379.4736
$ws.Range("I26").Value = 1531.7646
$ws.Range("J26").Value = 85
$ws.Range("K26").Value = 4595.293799999999
$ws.Range("L26").Value = 255
$ws.Range("M26").Value = -4307.293799999999
$ws.Range("N26").Value = -831
$ws.Range("H68").Value = 2177458.8
$ws.Range("J68").Value = 2503746
$ws.Range("L68").Value = 7511238
$ws.Range("N68").Value = -7512860
$ws.Range("H71").Value = 2177458.8
$ws.Range("J71").Value = 2503746
$ws.Range("L71").Value = 22533714
$ws.Range("N71").Value = -22541826
$ws.Range("H92").Value = 2500
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("H107").Value = 2817.102
$ws.Range("J107").Value = 3429.6052
$ws.Range("L107").Value = 10288.8156
$ws.Range("N107").Value = -14128.8156
$ws.Range("H114").Value = 2026.6875
$ws.Range("J114").Value = 2531.125
$ws.Range("L114").Value = 7593.375
$ws.Range("N114").Value = -14101.375
$ws.Range("H133").Value = 4727.75
$ws.Range("H134").Value = 2978.2856
$ws.Range("J134").Value = 5000
$ws.Range("L134").Value = 15000
$ws.Range("N134").Value = -25140
$ws.Range("H138").Value = 6034.3125
$ws.Range("I138").Value = 3134.2
$ws.Range("J138").Value = 10867.833
$ws.Range("K138").Value = 9402.599999999999
$ws.Range("L138").Value = 32603.499
$ws.Range("M138").Value = -4262.599999999999
$ws.Range("N138").Value = -42883.499

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 83399.664
$ws.Range("I2").Value = 83399.664
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 83399.664
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -83286.664
$ws.Range("N2").ClearContents()
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H6").Value = 1750
$ws.Range("I6").Value = 2500
$ws.Range("J6").Value = 1000
$ws.Range("K6").Value = 2500
$ws.Range("L6").Value = 1000
$ws.Range("M6").Value = -2387
$ws.Range("N6").Value = -1226
$ws.Range("H16").Value = 1750
$ws.Range("I16").Value = 2500
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 2500
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -2250
$ws.Range("N16").Value = -1500
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H70").Value = 125030000
$ws.Range("I70").Value = 500000000
$ws.Range("J70").Value = 40000
$ws.Range("K70").Value = 500000000
$ws.Range("L70").Value = 40000
$ws.Range("M70").Value = -499999730
$ws.Range("N70").Value = -40540
$ws.Range("H73").Value = 125030000
$ws.Range("I73").Value = 500000000
$ws.Range("J73").Value = 40000
$ws.Range("K73").Value = 500000000
$ws.Range("L73").Value = 40000
$ws.Range("M73").Value = -499999064
$ws.Range("N73").Value = -41872
$ws.Range("H97").Value = 838.6667
$ws.Range("I97").Value = 611.2105
$ws.Range("K97").Value = 611.2105
$ws.Range("M97").Value = -115.2105
$ws.Range("H102").Value = 111129210
$ws.Range("I102").Value = 125019110
$ws.Range("J102").Value = 10000
$ws.Range("K102").Value = 125019110
$ws.Range("L102").Value = 10000
$ws.Range("M102").Value = -125017488
$ws.Range("N102").Value = -13244
$ws.Range("H122").Value = 83339930
$ws.Range("I122").Value = 7192.1816
$ws.Range("J122").Value = 1000000000
$ws.Range("K122").Value = 21576.5448
$ws.Range("L122").Value = 3000000000
$ws.Range("M122").Value = -19126.5448
$ws.Range("N122").Value = -3000004900
$ws.Range("H126").Value = 83335330
$ws.Range("I126").Value = 100001800
$ws.Range("K126").Value = 300005400
$ws.Range("M126").Value = -300002930
$ws.Range("H132").Value = 4275
$ws.Range("I132").Value = 6000
$ws.Range("J132").Value = 2550
$ws.Range("K132").Value = 18000
$ws.Range("L132").Value = 7650
$ws.Range("M132").Value = -15470
$ws.Range("N132").Value = -12710
$ws.Range("H7").Value = 3000
$ws.Range("I7").Value = 3000
$ws.Range("K7").Value = 3000
$ws.Range("M7").Value = -2888

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 82419800
$ws.Range("I22").Value = 7938116
$ws.Range("J22").Value = 250003580
$ws.Range("K22").Value = 7938116
$ws.Range("L22").Value = 250003580
$ws.Range("M22").Value = -7937821
$ws.Range("N22").Value = -250004170
$ws.Range("H27").Value = 82419800
$ws.Range("I27").Value = 7938116
$ws.Range("J27").Value = 250003580
$ws.Range("K27").Value = 7938116
$ws.Range("L27").Value = 250003580
$ws.Range("M27").Value = -7938009
$ws.Range("N27").Value = -250003794
$ws.Range("H36").Value = 70000
$ws.Range("J36").Value = 70000
$ws.Range("L36").Value = 70000
$ws.Range("N36").Value = -71124
$ws.Range("H40").Value = 50380.5
$ws.Range("I40").Value = 82906.336
$ws.Range("J40").Value = 6027.091
$ws.Range("K40").Value = 82906.336
$ws.Range("L40").Value = 6027.091
$ws.Range("M40").Value = -82770.336
$ws.Range("N40").Value = -6299.091
$ws.Range("H68").Value = 2947.8235
$ws.Range("I68").Value = 2306.9167
$ws.Range("K68").Value = 2306.9167
$ws.Range("M68").Value = -1557.9167
$ws.Range("H71").Value = 2947.8235
$ws.Range("I71").Value = 2306.9167
$ws.Range("K71").Value = 11534.5835
$ws.Range("M71").Value = -7790.583500000001
$ws.Range("H82").Value = 1191.1666
$ws.Range("J82").Value = 1439.2858
$ws.Range("L82").Value = 1439.2858
$ws.Range("N82").Value = -2161.2858
$ws.Range("H85").Value = 1191.1666
$ws.Range("J85").Value = 1439.2858
$ws.Range("L85").Value = 1439.2858
$ws.Range("N85").Value = -3935.2858
$ws.Range("H100").Value = 3499
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("H122").Value = 5514
$ws.Range("I122").Value = 4617.2
$ws.Range("J122").Value = 9998
$ws.Range("K122").Value = 13851.6
$ws.Range("L122").Value = 29994
$ws.Range("M122").Value = -11401.6
$ws.Range("N122").Value = -34894
$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530
$ws.Range("H132").Value = 11237.28
$ws.Range("I132").Value = 11347.765
$ws.Range("J132").Value = 11002.5
$ws.Range("K132").Value = 34043.295
$ws.Range("L132").Value = 33007.5
$ws.Range("M132").Value = -31513.295
$ws.Range("N132").Value = -38067.5
$ws.Range("H136").Value = 5921.871
$ws.Range("I136").Value = 4983.36
$ws.Range("J136").Value = 9832.333000000001
$ws.Range("K136").Value = 14950.08
$ws.Range("L136").Value = 29496.999
$ws.Range("M136").Value = -12400.08
$ws.Range("N136").Value = -34596.999
$ws.Range("H139").Value = 66703.75
$ws.Range("J139").Value = 66703.75
$ws.Range("L139").Value = 66703.75
$ws.Range("N139").Value = -76983.75

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 32458.5
$ws.Range("J54").Value = 33688.25
$ws.Range("L54").Value = 33688.25
$ws.Range("N54").Value = -34728.25
$ws.Range("H64").Value = 891665.75
$ws.Range("J64").Value = 891665.75
$ws.Range("L64").Value = 891665.75
$ws.Range("N64").Value = -892161.75
$ws.Range("H67").Value = 891665.75
$ws.Range("J67").Value = 891665.75
$ws.Range("L67").Value = 891665.75
$ws.Range("N67").Value = -893381.75
$ws.Range("H107").Value = 679.1429000000001
$ws.Range("I107").Value = 679.1429000000001
$ws.Range("K107").Value = 2037.4287
$ws.Range("M107").Value = -117.4287000000002
$ws.Range("H122").Value = 9262427
$ws.Range("I122").Value = 3199.1667
$ws.Range("K122").Value = 9597.500100000001
$ws.Range("M122").Value = -7147.500100000001
$ws.Range("H132").Value = 5067.793
$ws.Range("I132").Value = 4961.815
$ws.Range("K132").Value = 14885.445
$ws.Range("M132").Value = -12355.445
